# Update odds values on Sheet1 per the 2024-10-14 FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7

# Row 3
$ws.Range("G3").Value = 2.05
$ws.Range("I3").Value = 4.33
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 8
$ws.Range("Z3").Value = 17
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 51
$ws.Range("AO3").Value = 12
$ws.Range("AU3").Value = 10
$ws.Range("AW3").Value = 6

# Row 4
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4

# Row 5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.5
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5

# Row 9
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.22
$ws.Range("P9").Value = 4
$ws.Range("Q9").Value = 1.8
$ws.Range("R9").Value = 2

$wb.Save()
